# Regenerate s_val data to filter save games: update computed stat columns
# for the two data rows (row 2 and row 3) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (2024-07-05)
$ws.Range("B2").Value = 1.459612070389937
$ws.Range("C2").Value = 25707020678.0705
$ws.Range("D2").Value = 26.21740644021617
$ws.Range("E2").Value = 8.660232485948974
$ws.Range("G2").Value = 25707020714.40776

# Row 3 (2024-06-21)
$ws.Range("B3").Value = 0.003994804209775715
$ws.Range("C3").Value = 0.04240448674262143
$ws.Range("D3").Value = 0.8054896365839992
$ws.Range("E3").Value = 645.3272768299601
$ws.Range("G3").Value = 646.1791657574964
